$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("C2").Value = 19
$ws.Range("C4").Value = 15.9

# Add new row 5 data
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2014-06-09"
$ws.Range("B5").Value = 123
$ws.Range("C5").Value = 22.5
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "hello "

# Update selection to match target state
$ws.Range("C5").Select()

$wb.Save()
